$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the "fazer login" (user + admin) tasks as done (column F = "FEITO") ---
# Rows 2, 3 (user login) and rows 12, 13 (admin login) get "REALIZADO",
# matching the formatting already used by e.g. F4/F6/F8 (green fill + centered).
$ws.Range("F2").Value = "REALIZADO"
$ws.Range("F3").Value = "REALIZADO"
$ws.Range("F12").Value = "REALIZADO"
$ws.Range("F13").Value = "REALIZADO"

$ws.Range("F4").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F13").PasteSpecial(-4122)

# --- Reassign the "editar perguntas" / "deletar perguntas" tasks to LILIAN ---
$ws.Range("E15").Value = "LILIAN"
$ws.Range("E16").Value = "LILIAN"

# --- Mark/select cell I26 (mirrors the existing empty marker cells G7/G8/I11) ---
$ws.Range("G7").Copy()
$ws.Range("I26").PasteSpecial(-4122)

$ws.Range("I26").Select()
